$d = $word.ActiveDocument

# 1. Capitalize "kit" -> "Kit" in "Loi Wei kit"
$d.Content.Find.Execute("Loi Wei kit", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Loi Wei Kit", 2)

# 2. Append ",LWK" after "LMJ"
$d.Content.Find.Execute("LMJ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "LMJ,LWK", 2)
